$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 2450
$ws.Range("I10").Value = 1900
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 1900
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -1607
$ws.Range("N10").Value = -3586

$ws.Range("H31").Value = 207.5
$ws.Range("I31").Value = 207.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 622.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -392.5

$ws.Range("H70").Value = 50584.8
$ws.Range("I70").Value = 1691.6666
$ws.Range("J70").Value = 246157.33
$ws.Range("K70").Value = 5074.9998
$ws.Range("L70").Value = 738471.99
$ws.Range("M70").Value = -4804.9998
$ws.Range("N70").Value = -739011.99

$ws.Range("H73").Value = 50584.8
$ws.Range("I73").Value = 1691.6666
$ws.Range("J73").Value = 246157.33
$ws.Range("K73").Value = 5074.9998
$ws.Range("L73").Value = 738471.99
$ws.Range("M73").Value = -4138.9998
$ws.Range("N73").Value = -740343.99

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H132").Value = 1375.75
$ws.Range("I132").Value = 1048.619
$ws.Range("J132").Value = 3665.6667
$ws.Range("K132").Value = 3145.857
$ws.Range("L132").Value = 10997.0001
$ws.Range("M132").Value = -615.857
$ws.Range("N132").Value = -16057.0001

$ws.Range("H138").Value = 4008.8772
$ws.Range("I138").Value = 3364.5938
$ws.Range("J138").Value = 4833.56
$ws.Range("K138").Value = 10093.7814
$ws.Range("L138").Value = 14500.68
$ws.Range("M138").Value = -4953.7814
$ws.Range("N138").Value = -24780.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1997.6666
$ws.Range("I61").Value = 1997.2
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1997.2
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1785.2
$ws.Range("N61").Value = -2424

$ws.Range("H122").Value = 669299.9
$ws.Range("I122").Value = 716714.1
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 2150142.3
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -2147692.3

$ws.Range("H132").Value = 964.8182
$ws.Range("I132").Value = 589.375
$ws.Range("J132").Value = 1966
$ws.Range("K132").Value = 1768.125
$ws.Range("L132").Value = 5898
$ws.Range("M132").Value = 761.875
$ws.Range("N132").Value = -10958

$ws.Range("H136").Value = 1997.6666
$ws.Range("I136").Value = 1997.2
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5991.6
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3441.6
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2019.875
$ws.Range("I20").Value = 1886
$ws.Range("J20").Value = 2421.5
$ws.Range("K20").Value = 1886
$ws.Range("L20").Value = 2421.5
$ws.Range("M20").Value = -1639

$ws.Range("H140").Value = 111111
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 111111
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 111111
$ws.Range("N140").Value = -121471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 395.15384
$ws.Range("I22").Value = 153
$ws.Range("J22").Value = 467.8
$ws.Range("K22").Value = 153
$ws.Range("L22").Value = 467.8
$ws.Range("M22").Value = 197

$ws.Range("H31").Value = 3528.25
$ws.Range("I31").Value = 3053.4
$ws.Range("J31").Value = 5902.5
$ws.Range("K31").Value = 3053.4
$ws.Range("L31").Value = 5902.5
$ws.Range("M31").Value = -2758.4

$ws.Range("H34").Value = 3528.25
$ws.Range("I34").Value = 3053.4
$ws.Range("J34").Value = 5902.5
$ws.Range("K34").Value = 3053.4
$ws.Range("L34").Value = 5902.5
$ws.Range("M34").Value = -2851.4

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("I132").Value = 1446.875
$ws.Range("J132").Value = 1150
$ws.Range("K132").Value = 4340.625
$ws.Range("L132").Value = 3450
$ws.Range("M132").Value = -1810.625
$ws.Range("N132").Value = -8510

$ws.Range("H134").Value = 1885.2162
$ws.Range("I134").Value = 1647.5483
$ws.Range("J134").Value = 3113.1667
$ws.Range("K134").Value = 4942.644899999999
$ws.Range("L134").Value = 9339.500100000001
$ws.Range("M134").Value = -2407.644899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 51.5
$ws.Range("I47").Value = 51.5
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 154.5
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 276.5

$ws.Range("H51").Value = 2460
$ws.Range("I51").Value = 2100
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 6300
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -5840

$ws.Range("H68").Value = 1484.375
$ws.Range("I68").Value = 1089.8
$ws.Range("J68").Value = 2142
$ws.Range("K68").Value = 3269.4
$ws.Range("L68").Value = 6426
$ws.Range("M68").Value = -2458.4
$ws.Range("N68").Value = -8048

$ws.Range("H71").Value = 1484.375
$ws.Range("I71").Value = 1089.8
$ws.Range("J71").Value = 2142
$ws.Range("K71").Value = 9808.199999999999
$ws.Range("L71").Value = 19278
$ws.Range("M71").Value = -5752.199999999999
$ws.Range("N71").Value = -27390

$ws.Range("H131").Value = 1331.3334
$ws.Range("I131").Value = 664.5
$ws.Range("J131").Value = 2665
$ws.Range("K131").Value = 1993.5
$ws.Range("L131").Value = 7995
$ws.Range("M131").Value = 3046.5
$ws.Range("N131").Value = -18075

$ws.Range("H137").Value = 4732.8335
$ws.Range("I137").Value = 3958.8
$ws.Range("J137").Value = 5285.7144
$ws.Range("K137").Value = 11876.4
$ws.Range("L137").Value = 15857.1432
$ws.Range("M137").Value = -6776.400000000001
$ws.Range("N137").Value = -26057.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6636.3184
$ws.Range("I70").Value = 6333.278
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 6333.278
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = -6063.278

$ws.Range("H73").Value = 6636.3184
$ws.Range("I73").Value = 6333.278
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 6333.278
$ws.Range("L73").Value = 8000
$ws.Range("M73").Value = -5397.278

$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1150
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1150
$ws.Range("M107").Value = 920

$ws.Range("H113").Value = 2283.3333
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 1300
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 870

$ws.Range("H132").Value = 2392.9656
$ws.Range("I132").Value = 2199.3333
$ws.Range("J132").Value = 2709.818
$ws.Range("K132").Value = 6597.999899999999
$ws.Range("L132").Value = 8129.454000000001
$ws.Range("M132").Value = -4067.999899999999
$ws.Range("N132").Value = -13189.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3563.56
$ws.Range("I22").Value = 2560.8635
$ws.Range("J22").Value = 10916.667
$ws.Range("K22").Value = 2560.8635
$ws.Range("L22").Value = 10916.667
$ws.Range("M22").Value = -2265.8635
$ws.Range("N22").Value = -11506.667

$ws.Range("H27").Value = 3563.56
$ws.Range("I27").Value = 2560.8635
$ws.Range("J27").Value = 10916.667
$ws.Range("K27").Value = 2560.8635
$ws.Range("L27").Value = 10916.667
$ws.Range("M27").Value = -2453.8635
$ws.Range("N27").Value = -11130.667

$ws.Range("H46").Value = 4686.5
$ws.Range("I46").Value = 2748.75
$ws.Range("J46").Value = 5655.375
$ws.Range("K46").Value = 2748.75
$ws.Range("L46").Value = 5655.375
$ws.Range("M46").Value = -2560.75
$ws.Range("N46").Value = -6031.375

$ws.Range("H58").Value = 450
$ws.Range("I58").Value = 450
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 450
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -190
$ws.Range("N58").ClearContents()

$ws.Range("H136").Value = 3074.35
$ws.Range("I136").Value = 2483.2778
$ws.Range("J136").Value = 8394
$ws.Range("K136").Value = 7449.8334
$ws.Range("L136").Value = 25182
$ws.Range("M136").Value = -4899.8334
$ws.Range("N136").Value = -30282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12931.417
$ws.Range("I81").Value = 13999.5
$ws.Range("J81").Value = 12717.8
$ws.Range("K81").Value = 27999
$ws.Range("L81").Value = 25435.6
$ws.Range("M81").Value = -26938
$ws.Range("N81").Value = -27557.6

$ws.Range("H84").Value = 12931.417
$ws.Range("I84").Value = 13999.5
$ws.Range("J84").Value = 12717.8
$ws.Range("K84").Value = 139995
$ws.Range("L84").Value = 127178
$ws.Range("M84").Value = -134691
$ws.Range("N84").Value = -137786

$ws.Range("H113").Value = 659.61536
$ws.Range("I113").Value = 540.4286
$ws.Range("J113").Value = 798.6667
$ws.Range("K113").Value = 1621.2858
$ws.Range("L113").Value = 2396.0001
$ws.Range("M113").Value = 548.7142000000001
